$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, pushing existing rows 8-15 down to 9-16.
$ws.Rows("8:8").Insert()

# Row 7 gets the corrected data (new shared strings are introduced here, in
# this order, so they land at the expected shared-string indexes).
$ws.Range("A7").Value = "Alex Correct Format"
$ws.Range("C7").Value = ".10/25/2003"
$ws.Range("D7").Value = "English"
$ws.Range("G7").Value = "1111/2222"

# Row 8 (new) becomes a copy of the former row 7 data, except column G now
# uses the new "1111/2222" string.
$ws.Range("A8").Value = "Alex Test6"
$ws.Range("B8").Value = "Male"
$ws.Range("C8").Value = ".10/28/2010"
$ws.Range("D8").Value = "Arabic"
$ws.Range("E8").Value = "Sudan"
$ws.Range("F8").Value = 1234567890
$ws.Range("G8").Value = "1111/2222"
$ws.Range("H8").Value = "RSD"

# Update the active selection to match the new layout.
$ws.Range("A6:H7").Select()
